# NYPD CompStat weekly report refresh:
#  - bump the "Volume/Number" edition text
#  - bump the reporting week date range
#  - refresh the crime-complaint data table (rows 15-31) with newly
#    collected figures, including a few cells that flip between a
#    numeric value and the "N/A" ("0") / "***.* " text placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: force a cell to hold a literal text value even when that text
# looks like a number (e.g. "0"), then restore the normal "text" look
# (font/alignment/General number format) by pasting formats from a
# cell that is already styled that way elsewhere in the sheet.
# ---------------------------------------------------------------------
function Set-TextSentinel($cellAddr, $text, $styleSourceAddr) {
    $ws.Range($cellAddr).NumberFormat = "@"
    $ws.Range($cellAddr).Value = $text
    $ws.Range($styleSourceAddr).Copy() | Out-Null
    $ws.Range($cellAddr).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Helper: turn a cell that currently holds placeholder text back into a
# real number, adopting the number format used by the rest of the
# numeric column (so the resulting style matches the other data cells).
# ---------------------------------------------------------------------
function Set-NumFromText($cellAddr, $value, $styleSourceAddr) {
    $fmt = $ws.Range($styleSourceAddr).NumberFormat
    $ws.Range($cellAddr).NumberFormat = $fmt
    $ws.Range($cellAddr).Value = $value
}

# ---------------------------------------------------------------------
# Header text updates
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/18/2024  Through  11/24/2024"

# ---------------------------------------------------------------------
# Row 15
# ---------------------------------------------------------------------
$ws.Range("N15").Value = -61.904761904761

# ---------------------------------------------------------------------
# Row 16
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 73.333333333333
$ws.Range("I16").Value = 224
$ws.Range("J16").Value = 147
$ws.Range("K16").Value = 52.380952380952
$ws.Range("L16").Value = 76.377952755905
$ws.Range("M16").Value = 39.130434782608
$ws.Range("N16").Value = -81.065088757396

# ---------------------------------------------------------------------
# Row 17
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 5.263157894736
$ws.Range("I17").Value = 304
$ws.Range("J17").Value = 238
$ws.Range("K17").Value = 27.731092436974
$ws.Range("L17").Value = 90
$ws.Range("M17").Value = 192.307692307692
$ws.Range("N17").Value = -23.809523809523

# ---------------------------------------------------------------------
# Row 18
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 127.272727272727
$ws.Range("I18").Value = 215
$ws.Range("J18").Value = 174
$ws.Range("K18").Value = 23.563218390804
$ws.Range("L18").Value = 29.518072289156
$ws.Range("M18").Value = 100.934579439252
$ws.Range("N18").Value = -71.179624664879

# ---------------------------------------------------------------------
# Row 19
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 220
$ws.Range("F19").Value = 77
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = 79.069767441860
$ws.Range("I19").Value = 719
$ws.Range("J19").Value = 621
$ws.Range("K19").Value = 15.780998389694
$ws.Range("L19").Value = 24.179620034542
$ws.Range("M19").Value = 53.961456102783
$ws.Range("N19").Value = -29.162561576354

# ---------------------------------------------------------------------
# Row 20 (C20 flips from a numeric value to the "N/A" text placeholder)
# ---------------------------------------------------------------------
Set-TextSentinel "C20" "0" "C14"
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 89
$ws.Range("J20").Value = 72
$ws.Range("K20").Value = 23.611111111111
$ws.Range("L20").Value = 48.333333333333
$ws.Range("M20").Value = 74.509803921568
$ws.Range("N20").Value = -86.873156342182

# ---------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 38.095238095238
$ws.Range("F21").Value = 156
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = 62.5
$ws.Range("I21").Value = 1561
$ws.Range("J21").Value = 1258
$ws.Range("K21").Value = 24.085850556438
$ws.Range("L21").Value = 41.523118766999
$ws.Range("M21").Value = 74.608501118568
$ws.Range("N21").Value = -61.485319516407

# ---------------------------------------------------------------------
# Row 22 (C22 flips num->text "N/A"; D22/E22 flip text->num)
# ---------------------------------------------------------------------
Set-TextSentinel "C22" "0" "C14"
Set-NumFromText "D22" 1 "F15"
Set-NumFromText "E22" -100 "H15"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 43
$ws.Range("J22").Value = 52
$ws.Range("K22").Value = -17.307692307692
$ws.Range("L22").Value = -23.214285714285
$ws.Range("M22").Value = -34.848484848484

# ---------------------------------------------------------------------
# Row 23 (C23 flips num->text "N/A"; D23/E23 flip text->num)
# ---------------------------------------------------------------------
Set-TextSentinel "C23" "0" "C14"
Set-NumFromText "D23" 1 "F15"
Set-NumFromText "E23" -100 "H15"
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 150
$ws.Range("I23").Value = 46
$ws.Range("J23").Value = 43
$ws.Range("K23").Value = 6.976744186046
$ws.Range("L23").Value = 31.428571428571
$ws.Range("M23").Value = 76.923076923076

# ---------------------------------------------------------------------
# Row 24
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 164
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = 24.242424242424
$ws.Range("I24").Value = 2012
$ws.Range("J24").Value = 1853
$ws.Range("K24").Value = 8.580679978413
$ws.Range("L24").Value = 15.300859598853
$ws.Range("M24").Value = 45.061283345349

# ---------------------------------------------------------------------
# Row 25
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 26
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 44.444444444444
$ws.Range("F25").Value = 143
$ws.Range("G25").Value = 105
$ws.Range("H25").Value = 36.190476190476
$ws.Range("I25").Value = 1790
$ws.Range("J25").Value = 1645
$ws.Range("K25").Value = 8.814589665653
$ws.Range("L25").Value = 23.193392980041

# ---------------------------------------------------------------------
# Row 26
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 114.285714285714
$ws.Range("F26").Value = 40
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = 5.263157894736
$ws.Range("I26").Value = 442
$ws.Range("J26").Value = 380
$ws.Range("K26").Value = 16.315789473684
$ws.Range("L26").Value = 32.732732732732
$ws.Range("M26").Value = 36

# ---------------------------------------------------------------------
# Row 27 (D27/E27 flip text->num)
# ---------------------------------------------------------------------
Set-NumFromText "D27" 1 "F15"
Set-NumFromText "E27" -100 "H15"
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 11
$ws.Range("K27").Value = 18.181818181818

# ---------------------------------------------------------------------
# Row 28 (C28 flips text->num; D28/E28 flip num->text "N/A"/"***.*")
# ---------------------------------------------------------------------
Set-NumFromText "C28" 3 "F15"
Set-TextSentinel "D28" "0" "C14"
Set-TextSentinel "E28" "***.*" "C14"
$ws.Range("I28").Value = 54
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -10

# ---------------------------------------------------------------------
# Row 29
# ---------------------------------------------------------------------
$ws.Range("N29").Value = -94.285714285714

# ---------------------------------------------------------------------
# Row 30
# ---------------------------------------------------------------------
$ws.Range("N30").Value = -92.307692307692

# ---------------------------------------------------------------------
# Row 31 (D31/E31 flip text->num; F31 flips num->text "N/A")
# ---------------------------------------------------------------------
Set-NumFromText "D31" 2 "F15"
Set-NumFromText "E31" -100 "H15"
Set-TextSentinel "F31" "0" "C14"
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = -100
$ws.Range("J31").Value = 20
$ws.Range("K31").Value = -10
